$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 437, pushing existing rows 437:463 down to 439:465
$ws.Rows("437:438").Insert()

# Row 437 - new weekly data point (Camote, 1a nueva(o), origin Peru)
$ws.Cells.Item(437, 1).Value = 11
$ws.Cells.Item(437, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(437, 3).Value = "Bíobío"
$ws.Cells.Item(437, 4).Value = 45223
$ws.Cells.Item(437, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(437, 5).Value = 8
$ws.Cells.Item(437, 6).Value = 100112045
$ws.Cells.Item(437, 7).Value = "Zapallo"
$ws.Cells.Item(437, 8).Value = "Camote"
$ws.Cells.Item(437, 9).Value = "1a nueva(o)"
$ws.Cells.Item(437, 10).Value = 500
$ws.Cells.Item(437, 11).Value = 1000
$ws.Cells.Item(437, 12).Value = 1000
$ws.Cells.Item(437, 13).Value = 1000
$ws.Cells.Item(437, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(437, 15).Value = "Perú"
$ws.Cells.Item(437, 16).Value = 1000
$ws.Cells.Item(437, 17).Value = 1
$ws.Cells.Item(437, 18).Value = "Hortaliza"

# Row 438 - new weekly data point (Paine, 1a (guarda), origin Region de O'Higgins)
$ws.Cells.Item(438, 1).Value = 11
$ws.Cells.Item(438, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(438, 3).Value = "Bíobío"
$ws.Cells.Item(438, 4).Value = 45223
$ws.Cells.Item(438, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(438, 5).Value = 8
$ws.Cells.Item(438, 6).Value = 100112045
$ws.Cells.Item(438, 7).Value = "Zapallo"
$ws.Cells.Item(438, 8).Value = "Paine"
$ws.Cells.Item(438, 9).Value = "1a (guarda)"
$ws.Cells.Item(438, 10).Value = 800
$ws.Cells.Item(438, 11).Value = 350
$ws.Cells.Item(438, 12).Value = 350
$ws.Cells.Item(438, 13).Value = 350
$ws.Cells.Item(438, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(438, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(438, 16).Value = 350
$ws.Cells.Item(438, 17).Value = 1
$ws.Cells.Item(438, 18).Value = "Hortaliza"
